$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet: the row-index column (A) is untouched; only the
#    date/count/value columns (B/C/D) shift down by one row, and a new row
#    9 is appended for the quarter that used to be the last one (2020-Q4).
#    A new "2022-Q3" entry is written into row 2.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Cells.Item(8, 1).Copy($total.Cells.Item(9, 1))
$total.Cells.Item(9, 1).Value = 7

for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    $total.Cells.Item($dest, 2).Value = $total.Cells.Item($r, 2).Value2
    $total.Cells.Item($dest, 3).Value = $total.Cells.Item($r, 3).Value2
    $total.Cells.Item($dest, 4).Value = $total.Cells.Item($r, 4).Value2
}

$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 43
$total.Cells.Item(2, 4).Value = 8.83

# ---------------------------------------------------------------------------
# 2. Insert a brand-new worksheet "2022-Q3" right after "总计", holding the
#    per-fund detail rows for the new quarter. Cloning the existing
#    "2022-Q2" sheet (instead of Worksheets.Add()) keeps the header /
#    column-A bold+border style ("s=2") intact without fabricating new
#    style entries.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# The cloned sheet only has 27 rows (26 data rows); the new quarter needs
# 44 (43 data rows), so extend column A (which carries the "s=2" style)
# down to row 44 by copying the existing styled A2 cell.
for ($r = 28; $r -le 44; $r++) {
    $q3.Cells.Item(2, 1).Copy($q3.Cells.Item($r, 1))
}

# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0,  "012541", "金鹰产业升级混合A",               "16.13", "93.57", "6.75", "1.0888", 6),
    @(1,  "519035", "富国天博创新混合",                 "18.16", "91.30", "5.02", "0.9116", 5),
    @(2,  "210008", "金鹰策略配置混合",                 "9.41",  "94.39", "9.14", "0.8601", 2),
    @(3,  "213001", "宝盈鸿利收益灵活配置混合A",         "16.06", "87.42", "4.91", "0.7885", 5),
    @(4,  "001298", "金鹰民族新兴灵活配置混合",          "8.10",  "93.37", "9.54", "0.7727", 2),
    @(5,  "007777", "中邮研究精选混合",                 "14.89", "70.15", "3.85", "0.5733", 8),
    @(6,  "002620", "中邮未来新蓝筹灵活配置混合",        "11.81", "83.48", "3.59", "0.4240", 9),
    @(7,  "000513", "富国高端制造行业股票A",             "7.01",  "91.41", "5.47", "0.3834", 5),
    @(8,  "210003", "金鹰行业优势混合",                 "5.59",  "82.52", "5.62", "0.3142", 2),
    @(9,  "011921", "富国均衡成长三年持有期混合A",       "6.14",  "90.69", "4.98", "0.3058", 2),
    @(10, "010751", "宝盈优质成长混合A",                "4.72",  "92.28", "5.97", "0.2818", 3),
    @(11, "008138", "富国龙头优势混合",                 "4.41",  "92.58", "5.23", "0.2306", 5),
    @(12, "014119", "金鹰时代先锋混合A",                "2.35",  "94.59", "9.70", "0.2280", 1),
    @(13, "008980", "中邮科技创新精选混合A",             "4.75",  "88.40", "4.58", "0.2176", 7),
    @(14, "590005", "中邮核心主题混合",                 "5.69",  "82.29", "3.58", "0.2037", 7),
    @(15, "673060", "西部利得景瑞灵活配置混合A",         "3.97",  "93.10", "4.34", "0.1723", 5),
    @(16, "001543", "宝盈新锐灵活配置混合A",            "2.20",  "91.14", "6.53", "0.1437", 3),
    @(17, "008981", "中邮科技创新精选混合C",            "2.66",  "88.40", "4.58", "0.1218", 7),
    @(18, "000264", "博时内需增长混合A",                "2.61",  "79.16", "3.65", "0.0953", 10),
    @(19, "160518", "博时睿远事件驱动灵活配置混合（LOF）", "2.64",  "78.13", "3.00", "0.0792", 6),
    @(20, "001277", "博时国企改革主题股票A",             "2.42",  "85.79", "3.08", "0.0745", 9),
    @(21, "050012", "博时策略混合",                    "2.18",  "79.24", "3.38", "0.0737", 10),
    @(22, "050014", "博时创业成长混合A",                "1.47",  "82.43", "4.11", "0.0604", 4),
    @(23, "009258", "西部利得景瑞灵活配置混合C",         "1.32",  "93.10", "4.34", "0.0573", 5),
    @(24, "005933", "新疆前海联合先进制造灵活配置混合A",  "0.95",  "92.14", "4.97", "0.0472", 4),
    @(25, "010752", "宝盈优质成长混合C",                "0.76",  "92.28", "5.97", "0.0454", 3),
    @(26, "012542", "金鹰产业升级混合C",                "0.63",  "93.57", "6.75", "0.0425", 6),
    @(27, "004265", "金鹰民丰回报定期开放混合",          "4.53",  "26.47", "0.78", "0.0353", 6),
    @(28, "007581", "宝盈鸿利收益灵活配置混合C",         "0.66",  "87.42", "4.91", "0.0324", 5),
    @(29, "160519", "博时睿利事件驱动灵活配置混合",      "0.55",  "82.58", "5.00", "0.0275", 3),
    @(30, "014120", "金鹰时代先锋混合C",                "0.26",  "94.59", "9.70", "0.0252", 1),
    @(31, "007578", "宝盈新锐灵活配置混合C",            "0.36",  "91.14", "6.53", "0.0235", 3),
    @(32, "011922", "富国均衡成长三年持有期混合C",       "0.44",  "90.69", "4.98", "0.0219", 2),
    @(33, "008061", "惠升惠新灵活配置混合A",            "0.40",  "82.68", "4.99", "0.0200", 5),
    @(34, "008533", "惠升惠兴混合A",                   "0.58",  "46.03", "3.18", "0.0184", 7),
    @(35, "004677", "博时战略新兴产业混合",             "0.38",  "91.02", "4.37", "0.0166", 6),
    @(36, "005934", "新疆前海联合先进制造灵活配置混合C",  "0.10",  "92.14", "4.97", "0.0050", 4),
    @(37, "008062", "惠升惠新灵活配置混合C",            "0.06",  "82.68", "4.99", "0.0030", 5),
    @(38, "002553", "博时创业成长混合C",                "0.05",  "82.43", "4.11", "0.0021", 4),
    @(39, "014930", "富国高端制造行业股票C",             "0.01",  "91.41", "5.47", "0.0005", 5),
    @(40, "011982", "博时内需增长混合C",                "0.01",  "79.16", "3.65", "0.0004", 10),
    @(41, "008534", "惠升惠兴混合C",                   "0.01",  "46.03", "3.18", "0.0003", 7),
    @(42, "014382", "博时国企改革主题股票C",             "0.00",  "85.79", "3.08", 0, 9)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    $r = $i + 2
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = "'" + $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    if ($i -eq ($rows.Length - 1)) {
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    }
    $q3.Cells.Item($r, 8).Value = $row[7]
}
